$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Find and delete the row whose keyword (column A) is "최대 당첨금"
$target = $ws.Columns.Item(1).Find("최대 당첨금")
if ($target -ne $null) {
    $target.EntireRow.Delete()
}

# Rename the keyword "복권 종류" to "복권" (its answer stays the same)
$target2 = $ws.Columns.Item(1).Find("복권 종류")
if ($target2 -ne $null) {
    $target2.Value = "복권"
}

# Find and delete the row whose keyword (column A) is "3배수"
$target3 = $ws.Columns.Item(1).Find("3배수")
if ($target3 -ne $null) {
    $target3.EntireRow.Delete()
}

# Update the frozen-pane / selection view state to match the saved workbook
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 8
$ws.Range("B27").Select()
